$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (shared strings content change)
$ws.Range("B1").Value = "ЕМБГ"
$ws.Range("C1").Value = "Датум на креирање"

# Remove the second data row entirely
$ws.Rows(2).Delete()

# Update column widths per target layout. Columns B and E revert to the
# workbook's default width; A, C and D get an explicit override.
$ws.Columns(1).ColumnWidth = 18.140625
$ws.Columns(2).ColumnWidth = 8.43
$ws.Columns(3).ColumnWidth = 19.85546875
$ws.Columns(4).ColumnWidth = 20.5703125
$ws.Columns(5).ColumnWidth = 8.43

# Update selection to match target sheet view
$ws.Range("E10").Select()
